$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers
$ws.Range("H1").Value = "A2-A3"
$ws.Range("I1").Value = "A2-A4"
$ws.Range("J1").Value = "A2-A5"
$ws.Range("K1").Value = "A2-A6"
$ws.Range("M1").Value = "A3-A4"
$ws.Range("N1").Value = "A3-A5"
$ws.Range("O1").Value = "A3-A6"
$ws.Range("Q1").Value = "A4-A5"
$ws.Range("R1").Value = "A4-A6"
$ws.Range("T1").Value = "A4-A6"

# Row 3 (Standard Deviation)
$ws.Range("H3").Value = 0.04377748576935908
$ws.Range("I3").Value = 0.1048904376256853
$ws.Range("J3").Value = 0.04465451977783355
$ws.Range("K3").Value = 0.1133286152801758
$ws.Range("M3").Value = 0.08199452009410384
$ws.Range("N3").Value = 0.01605221504347328
$ws.Range("O3").Value = 0.08615393724887715
$ws.Range("Q3").Value = 0.04363900915317774
$ws.Range("R3").Value = 0.03324744301279196
$ws.Range("T3").Value = 0.1096218056714369

# Row 4 (Maximum)
$ws.Range("H4").Value = 0.2302690318770986
$ws.Range("I4").Value = 0.4852792543902354
$ws.Range("J4").Value = 0.2183775035488023
$ws.Range("K4").Value = 0.4900937218377863
$ws.Range("M4").Value = 0.4454712637141095
$ws.Range("N4").Value = 0.07836700488541404
$ws.Range("O4").Value = 0.4509852702460233
$ws.Range("Q4").Value = 0.2414013973105197
$ws.Range("R4").Value = 0.1267255928289681
$ws.Range("T4").Value = 0.498967099440053

# Row 5 (Mean)
$ws.Range("A5").Value = "Mean"
$ws.Range("B5").Value = 0.01282128862737119
$ws.Range("C5").Value = 0.01554124594660452
$ws.Range("D5").Value = 0.05523844534376545
$ws.Range("E5").Value = 0.02773185877153201
$ws.Range("F5").Value = 0.06310082761723457
$ws.Range("H5").Value = 0.01942443474961239
$ws.Range("I5").Value = 0.05593333208081411
$ws.Range("J5").Value = 0.0240574250816833
$ws.Range("K5").Value = 0.055238142338576
$ws.Range("M5").Value = 0.03125602555220881
$ws.Range("N5").Value = 0.00967472014251436
$ws.Range("O5").Value = 0.03278094238855744
$ws.Range("Q5").Value = 0.01678073127719897
$ws.Range("R5").Value = 0.02183185037422592
$ws.Range("T5").Value = 0.05187175571702977

